# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header rows embedded as plain text in each sheet's A1/A2
# cell use PascalCase attribute names (ObjTablesVersion, Type, Id). This
# change lower-cases the first letter of those attribute names so they
# read as lowerCamelCase (objTablesVersion, type, id), matching the new
# convention. Only the textual content of these header cells changes;
# no structural changes are required.

$wb = $excel.ActiveWorkbook

# --- "!!Main root" sheet --------------------------------------------------
# Row 1 (A1): top-level "!!!ObjTables ObjTablesVersion='0.0.8'" header.
# Row 2 (A2): "!!ObjTables Type='Data' Id='MainRoot'" table header.
$wsMain = $wb.Worksheets.Item("!!Main root")
$wsMain.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsMain.Range("A2").Value = "!!ObjTables type='Data' id='MainRoot'"

# --- "!!Nodes" sheet -------------------------------------------------------
# Row 1 (A1): "!!ObjTables Type='Data' Id='Node'" table header.
$wsNodes = $wb.Worksheets.Item("!!Nodes")
$wsNodes.Range("A1").Value = "!!ObjTables type='Data' id='Node'"

# --- "!!Leaves" sheet -------------------------------------------------------
# Row 1 (A1): "!!ObjTables Type='Data' Id='Leaf'" table header.
$wsLeaves = $wb.Worksheets.Item("!!Leaves")
$wsLeaves.Range("A1").Value = "!!ObjTables type='Data' id='Leaf'"

# --- "!!One to many rows" sheet --------------------------------------------
# Row 1 (A1): "!!ObjTables Type='Data' Id='OneToManyRow'" table header.
$wsRows = $wb.Worksheets.Item("!!One to many rows")
$wsRows.Range("A1").Value = "!!ObjTables type='Data' id='OneToManyRow'"
